# assign4/data.xlsx edit: append "greedy" and "random" algorithm benchmark
# rows to Sheet1, rename the "Avg. moves" header, move the selection, and
# (best-effort) nudge the tab-ratio of the book view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- header rename: "Avg. moves (ms)" -> "Avg. moves" -------------------
$ws.Range("D1").Value = "Avg. moves"

# --- new benchmark rows ---------------------------------------------------
# greedy
$ws.Cells.Item(14, 1).Value = "greedy"
$ws.Cells.Item(14, 2).Value = 5
$ws.Cells.Item(14, 3).Value = 2.029
$ws.Cells.Item(14, 4).Value = 8.3333

$ws.Cells.Item(15, 1).Value = "greedy"
$ws.Cells.Item(15, 2).Value = 10
$ws.Cells.Item(15, 3).Value = 72.282
$ws.Cells.Item(15, 4).Value = 71.8

$ws.Cells.Item(16, 1).Value = "greedy"
$ws.Cells.Item(16, 2).Value = 15
$ws.Cells.Item(16, 3).Value = 213.587
$ws.Cells.Item(16, 4).Value = 83.3

$ws.Cells.Item(17, 1).Value = "greedy"
$ws.Cells.Item(17, 2).Value = 20
$ws.Cells.Item(17, 3).Value = 278.011
$ws.Cells.Item(17, 4).Value = 53.5

$ws.Cells.Item(18, 1).Value = "greedy"
$ws.Cells.Item(18, 2).Value = 25
$ws.Cells.Item(18, 3).Value = 751.173
$ws.Cells.Item(18, 4).Value = 82.8

$ws.Cells.Item(19, 1).Value = "greedy"
$ws.Cells.Item(19, 2).Value = 30
$ws.Cells.Item(19, 3).Value = 1209.045
$ws.Cells.Item(19, 4).Value = 83.4

$ws.Cells.Item(20, 1).Value = "greedy"
$ws.Cells.Item(20, 2).Value = 35
$ws.Cells.Item(20, 3).Value = 1969.064
$ws.Cells.Item(20, 4).Value = 90.7

$ws.Cells.Item(21, 1).Value = "greedy"
$ws.Cells.Item(21, 2).Value = 40
$ws.Cells.Item(21, 3).Value = 3435.527
$ws.Cells.Item(21, 4).Value = 110.6

$ws.Cells.Item(22, 1).Value = "greedy"
$ws.Cells.Item(22, 2).Value = 45
$ws.Cells.Item(22, 3).Value = 3874.541
$ws.Cells.Item(22, 4).Value = 91.6

$ws.Cells.Item(23, 1).Value = "greedy"
$ws.Cells.Item(23, 2).Value = 50
$ws.Cells.Item(23, 3).Value = 5889.61
$ws.Cells.Item(23, 4).Value = 104.2

$ws.Cells.Item(24, 1).Value = "greedy"
$ws.Cells.Item(24, 2).Value = 55
$ws.Cells.Item(24, 3).Value = 8539.096
$ws.Cells.Item(24, 4).Value = 115.9

$ws.Cells.Item(25, 1).Value = "greedy"
$ws.Cells.Item(25, 2).Value = 60
$ws.Cells.Item(25, 3).Value = 10362.746
$ws.Cells.Item(25, 4).Value = 111.4

# random
$ws.Cells.Item(26, 1).Value = "random"
$ws.Cells.Item(26, 2).Value = 5
$ws.Cells.Item(26, 3).Value = 2.597
$ws.Cells.Item(26, 4).Value = 11.4

$ws.Cells.Item(27, 1).Value = "random"
$ws.Cells.Item(27, 2).Value = 10
$ws.Cells.Item(27, 3).Value = 56.955
$ws.Cells.Item(27, 4).Value = 59.4

$ws.Cells.Item(28, 1).Value = "random"
$ws.Cells.Item(28, 2).Value = 15
$ws.Cells.Item(28, 3).Value = 284.025
$ws.Cells.Item(28, 4).Value = 118.5

$ws.Cells.Item(29, 1).Value = "random"
$ws.Cells.Item(29, 2).Value = 20
$ws.Cells.Item(29, 3).Value = 824.94
$ws.Cells.Item(29, 4).Value = 168

$ws.Cells.Item(30, 1).Value = "random"
$ws.Cells.Item(30, 2).Value = 25
$ws.Cells.Item(30, 3).Value = 1835.517
$ws.Cells.Item(30, 4).Value = 210.1

$ws.Cells.Item(31, 1).Value = "random"
$ws.Cells.Item(31, 2).Value = 30
$ws.Cells.Item(31, 3).Value = 2739.308
$ws.Cells.Item(31, 4).Value = 194.4

$ws.Cells.Item(32, 1).Value = "random"
$ws.Cells.Item(32, 2).Value = 35
$ws.Cells.Item(32, 3).Value = 3096.128
$ws.Cells.Item(32, 4).Value = 146.888

$ws.Cells.Item(33, 1).Value = "random"
$ws.Cells.Item(33, 2).Value = 40
$ws.Cells.Item(33, 3).Value = 5885.65
$ws.Cells.Item(33, 4).Value = 193.7

$ws.Cells.Item(34, 1).Value = "random"
$ws.Cells.Item(34, 2).Value = 45
$ws.Cells.Item(34, 3).Value = 10157.307
$ws.Cells.Item(34, 4).Value = 241.7

$ws.Cells.Item(35, 1).Value = "random"
$ws.Cells.Item(35, 2).Value = 50
$ws.Cells.Item(35, 3).Value = 11259.35
$ws.Cells.Item(35, 4).Value = 205.2

$ws.Cells.Item(36, 1).Value = "random"
$ws.Cells.Item(36, 2).Value = 55
$ws.Cells.Item(36, 3).Value = 13663.043
$ws.Cells.Item(36, 4).Value = 186.571

$ws.Cells.Item(37, 1).Value = "random"
$ws.Cells.Item(37, 2).Value = 60

# --- selection moves to F10 ------------------------------------------------
$ws.Range("F10").Select()

# --- book view tab ratio (best effort; 211 == 21.1%) -----------------------
$excel.ActiveWindow.TabRatio = 0.211
